$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text storage
# (source values like "6.90", "0.999" etc. are plain text in the sheet,
# not numbers -- force text format so Excel does not silently coerce them).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.013.08"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.913.86"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.79"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.89"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.90"
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000225"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.46"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.395.93"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.894.72"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.69"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.909.92"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.91"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.36"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.677"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.10"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.55"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.84"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.81"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.59"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.57"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.109"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0868"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.00"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.56"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.283"
$ws.Range("E41").Value = "  -4.54%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.94"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "378.09"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0346"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.697.26"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.72"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.79"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("E51").Value = "  -0.36%  "
